$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update employer number (column A) and member no (column B) for rows 2-7
for ($r = 2; $r -le 7; $r++) {
    $ws.Cells.Item($r, 1).Value = 10102368
    $ws.Cells.Item($r, 2).Value = 12345685
}

# Update the active selection to A2
$ws.Range("A2").Select()
